# Combine the "experienced vol" regression block with the master regression
# table: the exp_vol row moves to the top of the table (replacing rmse),
# the age-group rows now carry coefficients in every income/iqr column
# (B..G) instead of only D/G, and educ/HHinc rows shift down beneath them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper for cells whose text looks like a plain number (e.g. "40529",
# "0.02"). Excel would otherwise store these as numeric cells, but the
# source table keeps every value - including these - as text. Forcing the
# cell to Text format before assignment, then clearing the format again,
# keeps the value stored as a shared string without leaving a visible
# number-format change behind.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Helper for cells whose text is never ambiguous with a number (contains
# letters, "***", or parentheses) - these already round-trip as text.
function Set-TextValueStyled($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-TextValueStyled "B1" 'incvar I'
Set-TextValueStyled "C1" 'incvar II'
Set-TextValueStyled "D1" 'incvar III'
Set-TextValueStyled "E1" 'inciqr I'
Set-TextValueStyled "F1" 'inciqr II'
Set-TextValueStyled "G1" 'inciqr III'
Set-TextValueStyled "A2" 'exp_vol'
Set-TextValueStyled "B2" '1.81***'
Set-TextValueStyled "C2" '1.85***'
Set-TextValueStyled "D2" '1.95***'
Set-TextValueStyled "E2" '0.96***'
Set-TextValueStyled "F2" '0.94***'
Set-TextValueStyled "G2" '1.05***'
Set-TextValueStyled "B3" '(0.53)'
Set-TextValueStyled "C3" '(0.53)'
Set-TextValueStyled "D3" '(0.53)'
Set-TextValueStyled "E3" '(0.29)'
Set-TextValueStyled "F3" '(0.29)'
Set-TextValueStyled "G3" '(0.28)'
Set-TextValueStyled "A4" 'age_gr=30-39'
Set-TextValueStyled "B4" '-0.33***'
Set-TextValueStyled "C4" '-0.33***'
Set-TextValueStyled "D4" '-0.32***'
Set-TextValueStyled "E4" '-0.17***'
Set-TextValueStyled "F4" '-0.17***'
Set-TextValueStyled "G4" '-0.16***'
Set-TextValueStyled "B5" '(0.03)'
Set-TextValueStyled "C5" '(0.03)'
Set-TextValueStyled "D5" '(0.03)'
Set-TextValueStyled "E5" '(0.01)'
Set-TextValueStyled "F5" '(0.01)'
Set-TextValueStyled "G5" '(0.01)'
Set-TextValueStyled "A6" 'age_gr=40-48'
Set-TextValueStyled "B6" '-0.50***'
Set-TextValueStyled "C6" '-0.50***'
Set-TextValueStyled "D6" '-0.48***'
Set-TextValueStyled "E6" '-0.25***'
Set-TextValueStyled "F6" '-0.26***'
Set-TextValueStyled "G6" '-0.24***'
Set-TextValueStyled "B7" '(0.03)'
Set-TextValueStyled "C7" '(0.03)'
Set-TextValueStyled "D7" '(0.03)'
Set-TextValueStyled "E7" '(0.01)'
Set-TextValueStyled "F7" '(0.01)'
Set-TextValueStyled "G7" '(0.01)'
Set-TextValueStyled "A8" 'age_gr=49-57'
Set-TextValueStyled "B8" '-0.61***'
Set-TextValueStyled "C8" '-0.60***'
Set-TextValueStyled "D8" '-0.58***'
Set-TextValueStyled "E8" '-0.30***'
Set-TextValueStyled "F8" '-0.31***'
Set-TextValueStyled "G8" '-0.29***'
Set-TextValueStyled "B9" '(0.03)'
Set-TextValueStyled "C9" '(0.03)'
Set-TextValueStyled "D9" '(0.03)'
Set-TextValueStyled "E9" '(0.02)'
Set-TextValueStyled "F9" '(0.02)'
Set-TextValueStyled "G9" '(0.02)'
Set-TextValueStyled "A10" 'age_gr=>57'
Set-TextValueStyled "B10" '-0.48***'
Set-TextValueStyled "C10" '-0.47***'
Set-TextValueStyled "D10" '-0.45***'
Set-TextValueStyled "E10" '-0.23***'
Set-TextValueStyled "F10" '-0.25***'
Set-TextValueStyled "G10" '-0.23***'
Set-TextValueStyled "B11" '(0.04)'
Set-TextValueStyled "C11" '(0.04)'
Set-TextValueStyled "D11" '(0.04)'
Set-TextValueStyled "E11" '(0.02)'
Set-TextValueStyled "F11" '(0.02)'
Set-TextValueStyled "G11" '(0.02)'
Set-TextValueStyled "A12" 'educ_gr=low educ'
Set-TextValueStyled "C12" '-0.09***'
Set-TextValueStyled "D12" '-0.12***'
Set-TextValueStyled "F12" '0.07***'
Set-TextValueStyled "G12" '0.03***'
Set-TextValueStyled "C13" '(0.02)'
Set-TextValueStyled "D13" '(0.02)'
Set-TextValueStyled "F13" '(0.01)'
Set-TextValueStyled "G13" '(0.01)'
Set-TextValueStyled "A14" 'HHinc_gr=low inc'
Set-TextValueStyled "D14" '0.15***'
Set-TextValueStyled "G14" '0.19***'
Set-TextValueStyled "D15" '(0.02)'
Set-TextValueStyled "G15" '(0.01)'
Set-TextValueStyled "A16" 'N'
Set-TextValue "B16" '40529'
Set-TextValue "C16" '40529'
Set-TextValue "D16" '40529'
Set-TextValue "E16" '44874'
Set-TextValue "F16" '44874'
Set-TextValue "G16" '44874'
Set-TextValueStyled "A17" 'R2'
Set-TextValue "B17" '0.02'
Set-TextValue "C17" '0.02'
Set-TextValue "D17" '0.02'
Set-TextValue "E17" '0.02'
Set-TextValue "F17" '0.02'
Set-TextValue "G17" '0.03'
